$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.166.88"
$ws.Range("E2").Value = "  +0.91%  "
$ws.Range("D3").Value = "2.092.98"
$ws.Range("E3").Value = "  +2.86%  "
$ws.Range("E4").Value = "  -0.10%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "228.67"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.39%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.613"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.38%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "61.02"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +0.86%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +0.25%  "
$ws.Range("E10").Value = "  +3.75%  "
$ws.Range("E11").Value = "  -0.08%  "
$ws.Range("D12").Value = "2.402.96"
$ws.Range("E12").Value = "  +2.73%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "14.69"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +1.44%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "22.25"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +4.29%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "5.49"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +6.14%  "
$ws.Range("E16").Value = "  +2.02%  "
$ws.Range("D17").Value = "2.077.71"
$ws.Range("E17").Value = "  +1.75%  "
$ws.Range("D18").Value = "38.102.30"
$ws.Range("E18").Value = "  +0.68%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "6.02"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +2.11%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "70.17"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +0.54%  "
$ws.Range("E21").Value = "  +1.65%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "224.05"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("E24").Value = "  +0.57%  "
$ws.Range("E25").Value = "  +3.35%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "169.93"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +1.51%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "9.46"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +0.95%  "
$ws.Range("E28").Value = "  +0.15%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "18.96"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +0.49%  "
$ws.Range("E30").Value = "  +6.75%  "
$ws.Range("E31").Value = "  -0.26%  "
$ws.Range("E32").Value = "  +5.96%  "
$ws.Range("E33").Value = "  +4.33%  "
$ws.Range("E34").Value = "  +0.51%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.0606"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("E36").Value = "  +4.29%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "6.40"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +0.35%  "
$ws.Range("E38").Value = "  +4.65%  "
$ws.Range("E39").Value = "  -0.03%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "18.04"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +2.45%  "
$ws.Range("D41").Value = "1.554.85"
$ws.Range("E41").Value = "  +1.06%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "100.12"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +3.96%  "
$ws.Range("E43").Value = "  +0.37%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "2.84"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +0.95%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.0914"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +0.11%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "4.16"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +4.17%  "
$ws.Range("E47").Value = "  +1.33%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "7.44"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +4.76%  "
$ws.Range("E49").Value = "  +1.60%  "
$ws.Range("E50").Value = "  +0.81%  "
$ws.Range("D51").Value = "2.289.71"
$ws.Range("E51").Value = "  +2.74%  "
